$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the item rows that dropped out of this report run ---
# Delete from the bottom up so earlier row numbers stay valid while we work.
# Row 15: سرنجات 5 سم (item 9)
$ws.Rows.Item(15).EntireRow.Delete()
# Row 13: VOLTAREN 75MG/3ML 3 AMP. (item 7)
$ws.Rows.Item(13).EntireRow.Delete()
# Row 10: KAPRON 500 MG 20 F.C.TABS. (item 4)
$ws.Rows.Item(10).EntireRow.Delete()
# Row 9: DEXAMETHASONE-MUP 8MG/2ML 5 AMP (item 3)
$ws.Rows.Item(9).EntireRow.Delete()

# --- After the deletions the surviving items sit in rows 7-11; renumber them 1-5 ---
$ws.Range("A9").Value = 3
$ws.Range("A10").Value = 4
$ws.Range("A11").Value = 5

# --- Update the remaining سرنجات 3 سم row (now row 11) figures ---
$ws.Range("Q11").Value = "1:0"
# P11 holds a text-formatted number ("2.0000"); flip the cell to Text so the
# numeric-looking string isn't silently coerced into a real number, then
# restore the original number format so the style id is unchanged.
$origFmt = $ws.Range("P11").NumberFormat
$ws.Range("P11").NumberFormat = "@"
$ws.Range("P11").Value = "2.0000"
$ws.Range("P11").NumberFormat = $origFmt

# --- Update the totals row (now row 12) ---
$ws.Range("P12").Value = 89.74

# --- Update the generated timestamp in the footer row (now row 13) ---
$ws.Range("A13").Value = "Thursday, 21 August, 2025 10:45 AM"
